# Refresh cached market-price / profit figures on the Leve sheets
# (mirrors the scheduled market-data runner that produced the upstream commit)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 31
$ws.Range("H31").Value = 1979.8572
$ws.Range("I31").Value = 1979.8572
$ws.Range("K31").Value = 5939.571599999999
$ws.Range("M31").Value = -5709.571599999999

# Row 38
$ws.Range("H38").Value = 201.3077
$ws.Range("I38").Value = 209.75
$ws.Range("J38").Value = 100
$ws.Range("K38").Value = 629.25
$ws.Range("L38").Value = 300
$ws.Range("M38").Value = -257.25
$ws.Range("N38").Value = -1044

# Row 62
$ws.Range("H62").Value = 829.2857
$ws.Range("I62").Value = 884.1667
$ws.Range("J62").Value = 500
$ws.Range("K62").Value = 884.1667
$ws.Range("L62").Value = 500
$ws.Range("M62").Value = -260.1667
$ws.Range("N62").Value = -1748

# Row 65
$ws.Range("H65").Value = 829.2857
$ws.Range("I65").Value = 884.1667
$ws.Range("J65").Value = 500
$ws.Range("K65").Value = 4420.8335
$ws.Range("L65").Value = 2500
$ws.Range("M65").Value = -1300.8335
$ws.Range("N65").Value = -8740

# Row 76
$ws.Range("H76").Value = 10173.846
$ws.Range("I76").Value = 2626
$ws.Range("J76").Value = 35333.332
$ws.Range("K76").Value = 2626
$ws.Range("L76").Value = 35333.332
$ws.Range("M76").Value = -2311
$ws.Range("N76").Value = -35963.332

# Row 79
$ws.Range("H79").Value = 10173.846
$ws.Range("I79").Value = 2626
$ws.Range("J79").Value = 35333.332
$ws.Range("K79").Value = 2626
$ws.Range("L79").Value = 35333.332
$ws.Range("M79").Value = -1534
$ws.Range("N79").Value = -37517.332

# Row 125
$ws.Range("H125").Value = 5372.1055
$ws.Range("J125").Value = 5245.3335
$ws.Range("L125").Value = 47208.0015
$ws.Range("N125").Value = -52128.0015

# Row 137
$ws.Range("H137").Value = 971.83636
$ws.Range("I137").Value = 691.1429000000001
$ws.Range("J137").Value = 1012.7708
$ws.Range("K137").Value = 2073.4287
$ws.Range("L137").Value = 3038.3124
$ws.Range("M137").Value = 476.5712999999996
$ws.Range("N137").Value = -8138.3124

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 2268
$ws.Range("I63").Value = 2260
$ws.Range("J63").Value = 2380
$ws.Range("K63").Value = 2260
$ws.Range("L63").Value = 2380
$ws.Range("M63").Value = -1574
$ws.Range("N63").Value = -3752

# Row 66
$ws.Range("H66").Value = 2268
$ws.Range("I66").Value = 2260
$ws.Range("J66").Value = 2380
$ws.Range("K66").Value = 11300
$ws.Range("L66").Value = 11900
$ws.Range("M66").Value = -7868
$ws.Range("N66").Value = -18764

# Row 97
$ws.Range("H97").Value = 1506.0526
$ws.Range("I97").Value = 1240.2667
$ws.Range("J97").Value = 2502.75
$ws.Range("K97").Value = 1240.2667
$ws.Range("L97").Value = 2502.75
$ws.Range("M97").Value = -744.2666999999999
$ws.Range("N97").Value = -3494.75

$ws = $wb.Worksheets.Item("BSM")
# Row 56
$ws.Range("H56").Value = 68333.336
$ws.Range("J56").Value = 68333.336
$ws.Range("L56").Value = 68333.336
$ws.Range("N56").Value = -69811.336

# Row 86
$ws.Range("H86").Value = 1450
$ws.Range("I86").Value = 1440
$ws.Range("J86").Value = 1457.1428
$ws.Range("K86").Value = 1440
$ws.Range("L86").Value = 1457.1428
$ws.Range("M86").Value = -317
$ws.Range("N86").Value = -3703.1428

# Row 89
$ws.Range("H89").Value = 1450
$ws.Range("I89").Value = 1440
$ws.Range("J89").Value = 1457.1428
$ws.Range("K89").Value = 7200
$ws.Range("L89").Value = 7285.714
$ws.Range("M89").Value = -1584
$ws.Range("N89").Value = -18517.714

# Row 134
$ws.Range("H134").Value = 2553.2222
$ws.Range("I134").Value = 2471.8462
$ws.Range("J134").Value = 2764.8
$ws.Range("K134").Value = 7415.5386
$ws.Range("L134").Value = 8294.400000000001
$ws.Range("M134").Value = -4880.5386
$ws.Range("N134").Value = -13364.4

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 111111840
$ws.Range("I99").Value = 790
$ws.Range("J99").Value = 500000500
$ws.Range("K99").Value = 790
$ws.Range("L99").Value = 500000500
$ws.Range("M99").Value = 708
$ws.Range("N99").Value = -500003496

# Row 126
$ws.Range("H126").Value = 111111840
$ws.Range("I126").Value = 790
$ws.Range("J126").Value = 500000500
$ws.Range("K126").Value = 2370
$ws.Range("L126").Value = 1500001500
$ws.Range("M126").Value = 100
$ws.Range("N126").Value = -1500006440

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1290.52
$ws.Range("I68").Value = 978.14703
$ws.Range("J68").Value = 1549.561
$ws.Range("K68").Value = 2934.44109
$ws.Range("L68").Value = 4648.683
$ws.Range("M68").Value = -2123.44109
$ws.Range("N68").Value = -6270.683

# Row 71
$ws.Range("H71").Value = 1290.52
$ws.Range("I71").Value = 978.14703
$ws.Range("J71").Value = 1549.561
$ws.Range("K71").Value = 8803.323269999999
$ws.Range("L71").Value = 13946.049
$ws.Range("M71").Value = -4747.323269999999
$ws.Range("N71").Value = -22058.049

# Row 98
$ws.Range("H98").Value = 5052
$ws.Range("I98").Value = 100
$ws.Range("J98").Value = 10004
$ws.Range("K98").Value = 300
$ws.Range("L98").Value = 30012
$ws.Range("M98").Value = 1198
$ws.Range("N98").Value = -33008

# Row 121
$ws.Range("H121").Value = 29412628
$ws.Range("I121").Value = 400
$ws.Range("J121").Value = 31250892
$ws.Range("K121").Value = 1200
$ws.Range("L121").Value = 93752676
$ws.Range("M121").Value = 110
$ws.Range("N121").Value = -93755296

# Row 131
$ws.Range("H131").Value = 9616277
$ws.Range("I131").Value = 35714750
$ws.Range("J131").Value = 1048.6578
$ws.Range("K131").Value = 107144250
$ws.Range("L131").Value = 3145.9734
$ws.Range("M131").Value = -107139210
$ws.Range("N131").Value = -13225.9734

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2946.9285
$ws.Range("I122").Value = 2403.7693
$ws.Range("K122").Value = 7211.3079
$ws.Range("M122").Value = -4761.3079

# Row 126
$ws.Range("H126").Value = 1421.4117
$ws.Range("I126").Value = 1336.3334
$ws.Range("J126").Value = 1625.6
$ws.Range("K126").Value = 4009.0002
$ws.Range("L126").Value = 4876.799999999999
$ws.Range("M126").Value = -1539.0002
$ws.Range("N126").Value = -9816.799999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6465.1763
$ws.Range("I7").Value = 5646.769
$ws.Range("J7").Value = 9125
$ws.Range("K7").Value = 5646.769
$ws.Range("L7").Value = 9125
$ws.Range("M7").Value = -5534.769
$ws.Range("N7").Value = -9349

# Row 40
$ws.Range("H40").Value = 2722.1738
$ws.Range("I40").Value = 2567.4666
$ws.Range("J40").Value = 3012.25
$ws.Range("K40").Value = 2567.4666
$ws.Range("L40").Value = 3012.25
$ws.Range("M40").Value = -2431.4666
$ws.Range("N40").Value = -3284.25

# Row 46
$ws.Range("H46").Value = 1175.762
$ws.Range("I46").Value = 1194.2632
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 1194.2632
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -1006.2632
$ws.Range("N46").Value = -1376

# Row 126
$ws.Range("H126").Value = 6465.1763
$ws.Range("I126").Value = 5646.769
$ws.Range("J126").Value = 9125
$ws.Range("K126").Value = 16940.307
$ws.Range("L126").Value = 27375
$ws.Range("M126").Value = -14470.307
$ws.Range("N126").Value = -32315

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1883.5834
$ws.Range("I96").Value = 1850.5
$ws.Range("J96").Value = 1916.6666
$ws.Range("K96").Value = 1850.5
$ws.Range("L96").Value = 1916.6666
$ws.Range("M96").Value = -477.5
$ws.Range("N96").Value = -4662.6666

# Row 122
$ws.Range("H122").Value = 2221.1143
$ws.Range("I122").Value = 1960.3462
$ws.Range("K122").Value = 5881.0386
$ws.Range("M122").Value = -3431.0386

# Row 126
$ws.Range("H126").Value = 1783.0588
$ws.Range("I126").Value = 1782
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 5346
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -2876
$ws.Range("N126").Value = -10340

# Row 136
$ws.Range("H136").Value = 10914.083
$ws.Range("I136").Value = 3166
$ws.Range("J136").Value = 14788.125
$ws.Range("K136").Value = 9498
$ws.Range("L136").Value = 44364.375
$ws.Range("M136").Value = -6948
$ws.Range("N136").Value = -49464.375
